$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Swap the data (columns B..AC) between each of these row pairs. Column A
#    (the running index) is left untouched on both rows.
# ---------------------------------------------------------------------------
$pairs = @(
    @(14, 15),
    @(16, 17),
    @(64, 65),
    @(80, 81),
    @(87, 88),
    @(124, 125),
    @(137, 138),
    @(141, 142),
    @(156, 157),
    @(172, 173),
    @(186, 187),
    @(191, 192)
)

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]
    for ($c = 2; $c -le 29; $c++) {
        $v1 = $ws.Cells.Item($r1, $c).Value2
        $v2 = $ws.Cells.Item($r2, $c).Value2
        $ws.Cells.Item($r1, $c).Value2 = $v2
        $ws.Cells.Item($r2, $c).Value2 = $v1
    }
}

# ---------------------------------------------------------------------------
# 2) Rows 209 and 210 get results filled in (FTHG/FTAG/FTR) plus refreshed
#    odds-derived figures in columns R..AC.
# ---------------------------------------------------------------------------
$ws.Cells.Item(209, 8).Value2 = 1        # H209 (FTHG)
$ws.Cells.Item(209, 9).Value2 = 2        # I209 (FTAG)
$ws.Cells.Item(209, 10).Value2 = "A"     # J209 (FTR)
$ws.Cells.Item(209, 18).Value2 = 2.05    # R209
$ws.Cells.Item(209, 19).Value2 = 1.675   # S209
$ws.Cells.Item(209, 23).Value2 = -1      # W209
$ws.Cells.Item(209, 24).Value2 = -1      # X209
$ws.Cells.Item(209, 25).Value2 = 2.1     # Y209
$ws.Cells.Item(209, 26).Value2 = -1      # Z209
$ws.Cells.Item(209, 27).Value2 = 0.675   # AA209
$ws.Cells.Item(209, 28).Value2 = 0.825   # AB209
$ws.Cells.Item(209, 29).Value2 = -1      # AC209

$ws.Cells.Item(210, 8).Value2 = 2        # H210 (FTHG)
$ws.Cells.Item(210, 9).Value2 = 1        # I210 (FTAG)
$ws.Cells.Item(210, 10).Value2 = "H"     # J210 (FTR)
$ws.Cells.Item(210, 23).Value2 = 0.8     # W210
$ws.Cells.Item(210, 24).Value2 = -1      # X210
$ws.Cells.Item(210, 25).Value2 = -1      # Y210
$ws.Cells.Item(210, 26).Value2 = 0.825   # Z210
$ws.Cells.Item(210, 27).Value2 = -1      # AA210
$ws.Cells.Item(210, 28).Value2 = 1       # AB210
$ws.Cells.Item(210, 29).Value2 = -1      # AC210

# ---------------------------------------------------------------------------
# 3) Append two brand-new fixtures as rows 211 and 212 (not yet played, so no
#    FTHG/FTAG/FTR). Copy formats from row 210 so column A keeps the bold
#    centred "index" look and column E keeps the date/time number format.
# ---------------------------------------------------------------------------
$ws.Cells.Item(210, 1).Copy()
$ws.Cells.Item(211, 1).PasteSpecial(-4122)
$ws.Cells.Item(211, 1).Value2 = 209

$ws.Cells.Item(210, 5).Copy()
$ws.Cells.Item(211, 5).PasteSpecial(-4122)

$ws.Cells.Item(211, 2).Value2 = 7875109
$ws.Cells.Item(211, 3).Value2 = "Iraq League"
$ws.Cells.Item(211, 4).Value2 = "Iraq League"
$ws.Cells.Item(211, 5).Value2 = 45346.33333333334
$ws.Cells.Item(211, 6).Value2 = "Al Naft SC"
$ws.Cells.Item(211, 7).Value2 = "Al Minaa"
$ws.Cells.Item(211, 11).Value2 = 1.666
$ws.Cells.Item(211, 12).Value2 = 3.25
$ws.Cells.Item(211, 13).Value2 = 4.75
$ws.Cells.Item(211, 14).Value2 = 1.7
$ws.Cells.Item(211, 15).Value2 = 3.2
$ws.Cells.Item(211, 16).Value2 = 4.5
$ws.Cells.Item(211, 17).Value2 = -0.75
$ws.Cells.Item(211, 18).Value2 = 1.975
$ws.Cells.Item(211, 19).Value2 = 1.825
$ws.Cells.Item(211, 20).Value2 = 2.25
$ws.Cells.Item(211, 21).Value2 = 2
$ws.Cells.Item(211, 22).Value2 = 1.8
$ws.Cells.Item(211, 23).Value2 = 0
$ws.Cells.Item(211, 24).Value2 = 0
$ws.Cells.Item(211, 25).Value2 = 0
$ws.Cells.Item(211, 26).Value2 = 0
$ws.Cells.Item(211, 27).Value2 = 0

$ws.Cells.Item(210, 1).Copy()
$ws.Cells.Item(212, 1).PasteSpecial(-4122)
$ws.Cells.Item(212, 1).Value2 = 210

$ws.Cells.Item(210, 5).Copy()
$ws.Cells.Item(212, 5).PasteSpecial(-4122)

$ws.Cells.Item(212, 2).Value2 = 7875108
$ws.Cells.Item(212, 3).Value2 = "Iraq League"
$ws.Cells.Item(212, 4).Value2 = "Iraq League"
$ws.Cells.Item(212, 5).Value2 = 45346.4375
$ws.Cells.Item(212, 6).Value2 = "Al Zawraa"
$ws.Cells.Item(212, 7).Value2 = "Al Karkh"
$ws.Cells.Item(212, 11).Value2 = 1.615
$ws.Cells.Item(212, 12).Value2 = 3.25
$ws.Cells.Item(212, 13).Value2 = 5.25
$ws.Cells.Item(212, 14).Value2 = 1.666
$ws.Cells.Item(212, 15).Value2 = 3.2
$ws.Cells.Item(212, 16).Value2 = 5
$ws.Cells.Item(212, 17).Value2 = -0.75
$ws.Cells.Item(212, 18).Value2 = 1.95
$ws.Cells.Item(212, 19).Value2 = 1.85
$ws.Cells.Item(212, 20).Value2 = 1.75
$ws.Cells.Item(212, 21).Value2 = 1.925
$ws.Cells.Item(212, 22).Value2 = 1.875
$ws.Cells.Item(212, 23).Value2 = 0
$ws.Cells.Item(212, 24).Value2 = 0
$ws.Cells.Item(212, 25).Value2 = 0
$ws.Cells.Item(212, 26).Value2 = 0
$ws.Cells.Item(212, 27).Value2 = 0
